$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new "todo" rows (7 and 8) ---
# Copy the formatting of the prior data row (row 6) down onto the two new
# rows first, so the date cells (A/E) pick up the same built-in date
# number format (numFmtId 14) already used by the rest of the table.
$ws.Range("A6:E6").Copy()
$ws.Range("A7:E8").PasteSpecial(-4122)

# Row 7
$ws.Range("A7").Value = 42991
$ws.Range("B7").Value = "Improve code around calculating retentions"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 42991

# Row 8
$ws.Range("A8").Value = 42991
$ws.Range("B8").Value = "General code and function cleanup"
$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 6
$ws.Range("E8").Value = 42991

# --- Grow the worksheet table (Table1) so it covers the new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:E8"))

# --- Update the selection to match the author's final cursor position ---
[void]$ws.Range("B8").Select()
